$d = $word.ActiveDocument

# --- Edit 1: merge the " " run with the following
# "старший преподаватель Кленина Надежда Викторовна" run into a single run,
# without touching the preceding ", руководитель" run. ---
$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("старший преподаватель Кленина Надежда Викторовна")
$start1 = $rng1.Start - 1
$end1 = $rng1.End
$merge1 = $d.Range($start1, $end1)
$merge1.Delete()
$ins1 = $d.Range($start1, $start1)
$ins1.InsertAfter(" старший преподаватель Кленина Надежда Викторовна")

# --- Edit 2: merge the "в" run with the following
# " результате разработки было создано:" run into a single run,
# without touching the preceding bold/plain " " run. ---
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute(" результате разработки было создано:")
$start2 = $rng2.Start - 1
$end2 = $rng2.End
$merge2 = $d.Range($start2, $end2)
$merge2.Delete()
$ins2 = $d.Range($start2, $start2)
$ins2.InsertAfter("в результате разработки было создано:")

# --- Edit 3: insert a new bulleted paragraph "Алгоритм поиска пути"
# right after the "Генерация уровней" bullet, reusing that paragraph's
# list formatting. ---
$rng3 = $d.Content.Duplicate
$rng3.Find.Execute("Генерация уровней")
$rng3.Collapse(0)
$rng3.InsertParagraphAfter()

$rng3b = $d.Content.Duplicate
$rng3b.Find.Execute("Генерация уровней")
$genPara = $rng3b.Paragraphs(1)
$newPara = $genPara.Next()
$newPara.Range.InsertAfter("Алгоритм поиска пути")

# --- Edit 4: drop the stale "_GoBack" bookmark (Word relocates/clears
# this automatically as edits are made elsewhere in the document). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
